$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Re-use the existing "Landing Gear System Engineer" job description text
# (same value already stored in A3/A4) for the new row's A cell.
# Value2 is used (rather than Value) so the exact original string -
# including bullet characters - is copied verbatim and reuses the
# existing shared-string entry instead of minting a new one.
$jobDescription = $ws.Range("A3").Value2

$ws.Range("A6").Value2 = $jobDescription
$ws.Range("B6").Value = 1
$ws.Range("C6").Value = 7

# Entering a long multi-line string can make Excel stamp an explicit
# custom row height; AutoFit restores the default (no ht/customHeight
# attributes), matching the other data rows.
$ws.Rows.Item(6).EntireRow.AutoFit()
